$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric header labels (1..6) in B1:G1 with text labels G1..G6
$ws.Range("B1").Value = "G1"
$ws.Range("C1").Value = "G2"
$ws.Range("D1").Value = "G3"
$ws.Range("E1").Value = "G4"
$ws.Range("F1").Value = "G5"
$ws.Range("G1").Value = "G6"

# Update the selected cell to L3
$ws.Range("L3").Select()
